$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values remain plain text (they contain
# currency-group dots, not true decimal numbers) instead of Excel
# auto-converting them to numeric cells.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "26.570.25"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").Value = "1.820.91"
$ws.Range("E3").Value = "  +1.32%  "
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "1.006"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").Value = "307.53"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "0.4633"
$ws.Range("E7").Value = "  +2.09%  "
$ws.Range("D8").Value = "0.3598"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "0.07122"
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").Value = "0.8988"
$ws.Range("E10").Value = "  +1.26%  "
$ws.Range("D11").Value = "0.07762"
$ws.Range("E11").Value = "  -0.93%  "
$ws.Range("D12").Value = "19.35"
$ws.Range("E12").Value = "  -1.02%  "
$ws.Range("D13").Value = "1.826.80"
$ws.Range("E13").Value = "  +2.33%  "
$ws.Range("D14").Value = "5.249"
$ws.Range("E14").Value = "  -0.73%  "
$ws.Range("D15").Value = "6.307"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").Value = "87.22"
$ws.Range("E16").Value = "  +2.39%  "
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("D18").Value = "0.000008546"
$ws.Range("E18").Value = "  -0.47%  "
$ws.Range("D19").Value = "1.007"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").Value = "26.621.58"
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("D21").Value = "14.12"
$ws.Range("E21").Value = "  -1.33%  "
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("D24").Value = "1.913"
$ws.Range("E24").Value = "  -3.64%  "
$ws.Range("D25").Value = "151.92"
$ws.Range("E25").Value = "  -0.48%  "
$ws.Range("E26").Value = "  -0.37%  "
$ws.Range("D27").Value = "1.971"
$ws.Range("E27").Value = "  -4.06%  "
$ws.Range("D28").Value = "113.68"
$ws.Range("E28").Value = "  +1.31%  "
$ws.Range("D29").Value = "4.796"
$ws.Range("E29").Value = "  -1.68%  "
$ws.Range("D30").Value = "0.08799"
$ws.Range("E30").Value = "  +1.63%  "
$ws.Range("D31").Value = "3.135"
$ws.Range("E31").Value = "  +2.82%  "
$ws.Range("D32").Value = "0.7299"
$ws.Range("E32").Value = "  -0.07%  "
$ws.Range("D33").Value = "2.739"
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("D34").Value = "4.433"
$ws.Range("E34").Value = "  -0.67%  "
$ws.Range("D35").Value = "1.131"
$ws.Range("E35").Value = "  +1.55%  "
$ws.Range("D36").Value = "1.075"
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("E37").Value = "  -1.23%  "
$ws.Range("D38").Value = "2.915"
$ws.Range("E38").Value = "  +1.55%  "
$ws.Range("E39").Value = "  -0.44%  "
$ws.Range("D40").Value = "6.869"
$ws.Range("E40").Value = "  -0.24%  "
$ws.Range("D41").Value = "0.5043"
$ws.Range("E41").Value = "  -2.35%  "
$ws.Range("D42").Value = "0.1489"
$ws.Range("E42").Value = "  -2.13%  "
$ws.Range("D43").Value = "7.954"
$ws.Range("E43").Value = "  -0.78%  "
$ws.Range("D44").Value = "1.007"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").Value = "0.4642"
$ws.Range("E45").Value = "  -0.66%  "
$ws.Range("D46").Value = "9.924"
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("D47").Value = "98.06"
$ws.Range("E47").Value = "  -2.37%  "
$ws.Range("D48").Value = "1.555"
$ws.Range("E48").Value = "  -2.16%  "
$ws.Range("D49").Value = "0.05975"
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("D50").Value = "63.62"
$ws.Range("E50").Value = "  -1.21%  "
$ws.Range("E51").Value = "  -1.84%  "

# Restore the default cell style (the NumberFormat change above
# otherwise leaves the cells tagged with a non-default style index)
# while keeping the values stored as text.
$priceRange.Style = "Normal"
